# Sheet1's "String Prop" table (column G/H) lists pairs of label/value rows:
#   Younge's, Mass, Length, Radius ...
# This edit removes the "Length" row (H5:H6) entirely, so the "Radius"
# row that followed it (H7:H8) slides up to take its place. Doing this
# through direct value moves (rather than Range.Delete, which in this
# host shifts the *whole* row including column B) keeps column B intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value2 = $ws.Range("H7").Value2
$ws.Range("H6").Value2 = $ws.Range("H8").Value2
$ws.Range("H7").ClearContents()
$ws.Range("H8").ClearContents()

# Give column G a bit more breathing room and leave the selection there.
$ws.Columns("G").ColumnWidth = 10.67
$ws.Range("G10").Select()
